$d = $word.ActiveDocument

# --- Highlight 1: "Harvesting and forest floor removal..." ---
# "- Harvesting and forest floor removal alter arthropod, bacterial and fungal Communities in short term (5 year) time scales"
# -> "- Harvesting and forest floor removal alter soil communities in short term (5 year) time scales"
$d.Content.Find.Execute(
    "arthropod, bacterial and fungal Communities",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "soil communities",
    2) | Out-Null

# --- Highlight 4 (sampling location): ---
# "- Sampling location selection (relative to trenching) within sites can have a large effect on observed effects in sites where soil preparation has been applied "
# -> "- Sampling location within sites with soil preparation can have a large effect on observed effects "
$d.Content.Find.Execute(
    "Sampling location selection (relative to trenching) within sites can have a large effect on observed effects in sites where soil preparation has been applied",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Sampling location within sites with soil preparation can have a large effect on observed effects",
    2) | Out-Null

# --- Highlight 5 (DNA metabarcode): ---
# "- DNA metabarcode communities did not appear to be sensitive to seasonal fluctuations within the growing season, and may be a good technique to capture community patterns during the growing season regardless of the timing of sampling"
# -> "- Treatment responses in DNA metabarcode communities were consistent through different seasons"
# Split into two replacements so the proofErr-wrapped "metabarcode" run is left untouched.
$d.Content.Find.Execute(
    "- DNA ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "- Treatment responses in DNA ",
    2) | Out-Null

$d.Content.Find.Execute(
    " communities did not appear to be sensitive to seasonal fluctuations within the growing season, and may be a good technique to capture community patterns during the growing season regardless of the timing of sampling",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " communities were consistent through different seasons",
    2) | Out-Null
